$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6162.199529381775
$ws.Range("C2").Value = 2808.369209570993
$ws.Range("D2").Value = 6023.487960653602
